$d = $word.ActiveDocument
$replacements = @(
    @("2023-06-17 Saturday", "2023-06-18 Sunday"),
    @("75-26=49", "6+40=46"),
    @("75-13=62", "33+55=88"),
    @("54+24=78", "22+18=40"),
    @("89-73=16", "60-14=46"),
    @("11-2=9", "16+37=53"),
    @("58+11=69", "34-31=3"),
    @("44-9=35", "80-62=18"),
    @("29-20=9", "50+18=68"),
    @("47+46=93", "50+8=58"),
    @("30+7=37", "31-15=16"),
    @("70-50=20", "87+5=92"),
    @("32+35=67", "35+51=86"),
    @("45-22=23", "44-18=26"),
    @("28+3=31", "14+12=26"),
    @("67+28=95", "4+93=97"),
    @("0+48=48", "6+61=67"),
    @("73-71=2", "92-22=70"),
    @("58+4=62", "36-10=26"),
    @("62-59=3", "93-67=26"),
    @("4+89=93", "31+65=96"),
    @("17+49=66", "66-6=60"),
    @("81-15=66", "47-33=14"),
    @("76-64=12", "80-14=66"),
    @("46+4=50", "81-35=46"),
    @("93-10=83", "44-42=2"),
    @("8+12=20", "59-30=29"),
    @("44+14=58", "88-5=83"),
    @("65+8=73", "37+19=56"),
    @("12+36=48", "82-71=11"),
    @("44-40=4", "19+28=47"),
    @("82-5=77", "68-30=38"),
    @("20-6=14", "27-5=22"),
    @("64-43=21", "49-42=7"),
    @("69+13=82", "58-56=2"),
    @("41-16=25", "47+47=94"),
    @("78+11=89", "47-40=7"),
    @("30+47=77", "40+21=61"),
    @("87-13=74", "7-1=6"),
    @("23+4=27", "10+82=92"),
    @("34-13=21", "48+50=98"),
    @("77-58=19", "68-24=44"),
    @("90-67=23", "81-5=76"),
    @("27-26=1", "56+17=73"),
    @("0+11=11", "45-40=5"),
    @("39+21=60", "16+61=77"),
    @("50+29=79", "85+4=89"),
    @("38-35=3", "71-69=2"),
    @("39+26=65", "37+32=69"),
    @("67-25=42", "55-38=17"),
    @("21-7=14", "89-61=28"),
    @("34+6=40", "49-4=45"),
    @("16-2=14", "75-8=67"),
    @("70+5=75", "76-48=28"),
    @("64+26=90", "45-24=21"),
    @("57+17=74", "0+24=24"),
    @("25+23=48", "9+20=29"),
    @("90+3=93", "68-9=59"),
    @("8+83=91", "26+20=46"),
    @("44-14=30", "86-17=69"),
    @("62-42=20", "40+37=77"),
    @("83-37=46", "88-24=64"),
    @("44-10=34", "72+2=74"),
    @("61-28=33", "67-24=43"),
    @("28+68=96", "76-13=63"),
    @("15+33=48", "26-23=3"),
    @("17+58=75", "5-3=2"),
    @("76-32=44", "30+36=66"),
    @("22+63=85", "32-29=3"),
    @("35+56=91", "94-47=47"),
    @("84-14=70", "85-71=14"),
    @("99-84=15", "9+39=48"),
    @("31+62=93", "84-36=48"),
    @("38-32=6", "28+56=84"),
    @("20+71=91", "63-15=48"),
    @("82-7=75", "14+3=17"),
    @("37-18=19", "2+97=99"),
    @("59+11=70", "40-19=21"),
    @("99-38=61", "65+27=92"),
    @("28+27=55", "31+44=75"),
    @("49+32=81", "9+23=32"),
    @("62+10=72", "52+19=71"),
    @("6+42=48", "3+13=16"),
    @("53-29=24", "98-74=24"),
    @("37-17=20", "17+51=68"),
    @("54-35=19", "82-11=71"),
    @("94-59=35", "78-3=75"),
    @("48-23=25", "28+14=42"),
    @("66+31=97", "87-30=57"),
    @("47+40=87", "38+6=44"),
    @("68+15=83", "46-36=10"),
    @("27+9=36", "21+46=67"),
    @("17+63=80", "0+87=87"),
    @("65+21=86", "30+28=58"),
    @("67-42=25", "97-46=51"),
    @("30+1=31", "84-43=41"),
    @("65+25=90", "12-7=5"),
    @("16+71=87", "21+72=93"),
    @("72-53=19", "76-29=47"),
    @("79-1=78", "86-35=51"),
    @("11+41=52", "73-18=55"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Output "Replaced $($replacements.Count) items"